$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.232.25"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  +0.26%  "

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.70"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  +0.15%  "

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.03%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6988"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.81%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.80"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("E7").Value = "  +0.07%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08120"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +5.04%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3017"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -0.77%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.47"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +1.36%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08180"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.27%  "

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.856.30"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -0.76%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.193"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -0.01%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7057"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -2.33%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.73"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +0.43%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.249.59"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +0.38%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.818"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +1.61%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007903"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +1.19%  "

$ws.Range("E19").Value = "  +0.70%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.99"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +1.35%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -0.03%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.107.79"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.41%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -0.02%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.442"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.01%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.89"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.78%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.874"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.85%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1415"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -1.06%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.06"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +0.08%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.915"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -2.37%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.411"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +0.67%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.475"
$ws.Range("D31").Style = $origStyle

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.355"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -3.74%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.027"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("E34").Value = "  -0.14%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.161"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("E36").Value = "  +2.43%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9977"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -2.73%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.687"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +1.07%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01849"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +0.08%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.715"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +1.35%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9331"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +1.57%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.146.60"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +4.04%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.993"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -0.05%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4250"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("E46").Value = "  +0.06%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.86"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +0.49%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5287"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -3.60%  "

$ws.Range("E49").Value = "  -0.27%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.008.23"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +0.79%  "

$ws.Range("E51").Value = "  -0.06%  "

